# Adds three new "0x87" Heading1 + Given/When/Then/Additional table
# sections to the end of the document, right before the trailing empty
# Heading1 paragraph that closes the body.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Escape-Xml($text) {
    $t = $text -replace '&', '&amp;'
    $t = $t -replace '<', '&lt;'
    $t = $t -replace '>', '&gt;'
    return $t
}

function Build-RunXml($text) {
    if ($text -eq "") {
        return ""
    }
    $escaped = Escape-Xml $text
    $needsPreserve = ($text -ne $text.Trim()) -or ($text.Contains("  "))
    if ($needsPreserve) {
        return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    } else {
        return "<w:r><w:t>$escaped</w:t></w:r>"
    }
}

function Build-CellXml($text, $width) {
    $runXml = Build-RunXml $text
    if ($runXml -eq "") {
        $pXml = "<w:p/>"
    } else {
        $pXml = "<w:p>$runXml</w:p>"
    }
    return "<w:tc><w:tcPr><w:tcW w:w=`"$width`" w:type=`"dxa`"/></w:tcPr>$pXml</w:tc>"
}

function Build-TableXml($rows) {
    $trXml = ""
    foreach ($row in $rows) {
        $left = Build-CellXml $row[0] 1194
        $right = Build-CellXml $row[1] 7822
        $trXml += "<w:tr>$left$right</w:tr>"
    }
    $tblPr = '<w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/>' +
             '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr>'
    $tblGrid = '<w:tblGrid><w:gridCol w:w="1194"/><w:gridCol w:w="7822"/></w:tblGrid>'
    return "<w:tbl>$tblPr$tblGrid$trXml</w:tbl>"
}

function Get-EndInsertionRange() {
    $lastPara = $d.Paragraphs.Last
    return $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
}

function Insert-BodyXml($fragmentXml) {
    $ip = Get-EndInsertionRange
    $xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>$fragmentXml</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
    [void]$ip.InsertXML($xmlFrag)
}

function Add-TestCaseSection($heading, $rows) {
    # New Heading1 paragraph, inserted as a genuine paragraph break before
    # the trailing (final) paragraph of the document.
    $ip = Get-EndInsertionRange
    $ip.InsertParagraphBefore()
    $newHeading = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $newHeading.Range.ParagraphFormat.Style = "Heading1"
    $newHeading.Range.InsertAfter($heading)

    # New Given/When/Then/Additional table, inserted via raw OOXML so the
    # table properties match the existing tables exactly.
    $tblXml = Build-TableXml $rows
    Insert-BodyXml $tblXml
}

# ---------------------------------------------------------------------
# 0x87: Test_RDGN_INVALID_DIAG
# ---------------------------------------------------------------------

Add-TestCaseSection "0x87: Test_RDGN_INVALID_DIAG" @(
    @("Given:", "The node number of the unit under test, a valid service index number and an invalid diagnostic number"),
    @("When:", "The cbus message RDGN is sent"),
    @("Then:", "Expect a GRSP message with a result of ‘Invalid Diagnostic’ (253)"),
    @("Additional:", "")
)

# ---------------------------------------------------------------------
# 0x87: Test_RDGN_INVALID_SERVICE
# ---------------------------------------------------------------------

Add-TestCaseSection "0x87: Test_RDGN_INVALID_SERVICE" @(
    @("Given:", "The node number of the unit under test, an invalid service index number and an valid diagnostic number"),
    @("When:", "The cbus message RDGN is sent"),
    @("Then:", "Expect a GRSP message with a result of ‘Invalid Service’ (252)"),
    @("Additional:", "")
)

# ---------------------------------------------------------------------
# 0x87: Test_RDGN_SHORT
# ---------------------------------------------------------------------

Add-TestCaseSection "0x87: Test_RDGN_SHORT" @(
    @("Given:", "The node number of the unit under test, a valid service index number and a valid diagnostic number"),
    @("When:", "The cbus message RDGN is sent with the diagnostic number missing"),
    @("Then:", "Expect a GRSP message with a result of ‘Invalid Command’ (?)"),
    @("Additional:", "")
)

Write-Host "Done. Tables:" $d.Tables.Count "Paragraphs:" $d.Paragraphs.Count
